# Generate Report for Handoff
# Marks the "fc4e3003-cc78-4da1-a49e-947c16d8d7b7" file's status as
# "Ready for handoff" (was "In Translation") on all sheets, and updates
# the related handoff datetime stamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet, row 3 = fc4e3003-cc78-4da1-a49e-947c16d8d7b7.md
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-23 18:19:11"

# zh-cn detail sheet, row 3 = fc4e3003-cc78-4da1-a49e-947c16d8d7b7.md
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-23 18:19:07"

# de-de detail sheet, row 3 = fc4e3003-cc78-4da1-a49e-947c16d8d7b7.md
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-23 18:19:11"
